# Updates cryptos list values per commit "Updated cryptos list on Wed Sep 11 04:52:43 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.359.51"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "'2.325.44"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'511.19"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "'132.11"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").Value = "'0.0999"
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").Value = "'5.23"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "'0.336"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "'2.741.40"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "'23.50"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "'56.338.65"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").Value = "'2.329.08"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'10.34"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'322.67"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "'4.14"
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").Value = "'6.58"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'61.14"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("B24").Value = "Kaspa"
$ws.Range("C24").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D24").Value = "'0.163"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'8.53"
$ws.Range("E25").Value = "  +10.49%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "'1.29"
$ws.Range("E27").Value = "  +5.99%  "
$ws.Range("D28").Value = "'167.09"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.66"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0₃0714"
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").Value = "'6.08"
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("D32").Value = "'18.23"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("D36").Value = "'3.93"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").Value = "'0.879"
$ws.Range("E37").Value = "  -5.87%  "
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").Value = "'149.88"
$ws.Range("E40").Value = "  +8.79%  "
$ws.Range("D41").Value = "'0.373"
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("D42").Value = "'3.53"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "'276.52"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "'5.03"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").Value = "'0.0923"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").Value = "'0.0494"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").Value = "'0.552"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "'17.99"
$ws.Range("E48").Value = "  +5.53%  "
$ws.Range("D49").Value = "'0.380"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").Value = "'16.87"
$ws.Range("E51").Value = "  +1.14%  "
